# Update "想去人数" (column F) counts to the freshly scraped values.
# Sheet name -> list of (cell, newValue) pairs.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ "F2" = 1204; "F5" = 1346; "F6" = 1705; "F7" = 6226; "F8" = 126; "F9" = 1820; "F15" = 25; "F16" = 6943; "F21" = 1708; "F26" = 1584; "F27" = 761 }
    "演出"   = @{ "F4" = 348 }
    "本地生活" = @{ "F2" = 9531; "F3" = 2262; "F4" = 659; "F5" = 246 }
    "全部类型" = @{ "F2" = 9531; "F3" = 2262; "F4" = 659; "F5" = 1204; "F9" = 348; "F10" = 1346; "F11" = 246; "F12" = 1705; "F13" = 6226; "F14" = 1820; "F23" = 6943; "F28" = 1708; "F33" = 1584; "F34" = 761 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $updates[$sheetName]
    foreach ($cellRef in $cells.Keys) {
        $ws.Range($cellRef).Value = $cells[$cellRef]
    }
}

$wb.Save()
